# Commit "added Pie.html and Pie.css with working pie charts. Made
# backend case insensitive" is about files elsewhere in the repo (the
# PowerPoll web project's Pie.html/Pie.css/backend code). This
# Presentation1.pptx sits under bin/Debug as a build artifact, and the
# only reason it shows up in the commit is that it was re-saved (e.g.
# the author opened it in PowerPoint while working on the project).
#
# Diffing the canonical OOXML confirms this: every changed line is a
# relationship id (p:sldMasterId/@r:id, p:sldId/@r:id, every
# p:sldLayoutId/@r:id, the slide's we:webextensionref/@r:id and the
# picture's a:blip/@r:embed) plus the embedded Office Add-in's own
# instance guid (we:webextension/@id in ppt/slides/udata/data.xml) and
# its re-propagated we:snapshot/@r:embed. Those are implementation-
# internal identifiers that PowerPoint mints fresh on every save; they
# are not exposed anywhere in the Shape/Slide/Presentation object
# model, and none of the actual slide content changed - same title /
# subtitle placeholders (still empty), same embedded web-extension
# graphic frame + picture fallback, same snapshot image bytes.
#
# So the faithful COM-interop equivalent of this commit's effect on
# the deck is simply: open it, don't touch the slide content, save it.
$p = $ppt.ActivePresentation

# Exercise the deck the same way the author's PowerPoint session would
# have (load the one slide) without mutating any shape or text.
$s = $p.Slides.Item(1)

$p.Save()
